$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Cells.Item(1, 3).Value = "Color"

# Policy rows (2-21) -> #8FBC8F
for ($i = 2; $i -le 21; $i++) {
    $ws.Cells.Item($i, 3).Value = "#8FBC8F"
}

# Barrier rows (22-46) -> #B22222
for ($i = 22; $i -le 46; $i++) {
    $ws.Cells.Item($i, 3).Value = "#B22222"
}

# Concern rows (47-64) -> #DAA520
for ($i = 47; $i -le 64; $i++) {
    $ws.Cells.Item($i, 3).Value = "#DAA520"
}

# Update the active selection to match the saved workbook's view state
[void]$ws.Range("I62").Select()
